# Loan RBI, Variable Instalments
#
# 1. On the "Repayment Schedule" sheet, insert a new blank column before
#    column N (pushes the old N/O/P "Late" / "Heading" / "Outstanding"
#    columns one slot to the right, to O/P/Q) -- adds a "Variable
#    Instalment" style gap column to the repayment schedule.
# 2. Make "Repayment Schedule" the active sheet/tab (was "Input").
# 3. Update the remembered selection + scroll position on "Repayment
#    Schedule" to T8 (was I5), scrolled so column F is visible.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N (existing N/O/P shift right to O/P/Q).
$wsSchedule.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet (this both sets
# workbook.xml's bookViews/workbookView activeTab, and flips
# sheetView/tabSelected off "Input" and onto "Repayment Schedule").
$wsSchedule.Activate()

# Scroll so column F is the left-most visible column (was F, stays F,
# but row resets to the top) and move the selection to T8.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$wsSchedule.Range("T8").Select()
